$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (Price/Volume columns refreshed by the scraper).
$updates = @{
    'D2' = '25.862.84'
    'D3' = '1.736.08'
    'E3' = '  -0.39%  '
    'E4' = '  +0.00%  '
    'D5' = '242.38'
    'E5' = '  +5.11%  '
    'D6' = '1.000'
    'E6' = '  +0.06%  '
    'D7' = '0.5200'
    'E7' = '  -1.06%  '
    'D8' = '0.2740'
    'E8' = '  -0.64%  '
    'D9' = '0.06160'
    'E9' = '  +0.15%  '
    'D10' = '1.739.23'
    'E10' = '  -0.17%  '
    'D11' = '0.07187'
    'E11' = '  +1.29%  '
    'D12' = '14.97'
    'E12' = '  -2.08%  '
    'D13' = '0.6433'
    'E13' = '  -0.07%  '
    'D14' = '4.612'
    'E14' = '  +1.71%  '
    'D15' = '77.29'
    'E15' = '  -0.37%  '
    'D16' = '1.000'
    'E16' = '  +0.04%  '
    'D17' = '0.9998'
    'E17' = '  +0.00%  '
    'D18' = '25.896.92'
    'D19' = '11.76'
    'E19' = '  +1.93%  '
    'D20' = '0.000006780'
    'E20' = '  +1.32%  '
    'D21' = '1.962.40'
    'E21' = '  -0.05%  '
    'D22' = '4.283'
    'E22' = '  -0.48%  '
    'D23' = '8.632'
    'E23' = '  -1.74%  '
    'D24' = '5.290'
    'E24' = '  +2.44%  '
    'D25' = '137.16'
    'E25' = '  -2.17%  '
    'D26' = '1.510'
    'E26' = '  -0.70%  '
    'E27' = '  +0.36%  '
    'D28' = '1.773'
    'E28' = '  -1.33%  '
    'D29' = '105.19'
    'E29' = '  +2.42%  '
    'D30' = '3.968'
    'E30' = '  +6.45%  '
    'D31' = '0.08253'
    'E31' = '  -0.71%  '
    'D32' = '3.649'
    'E32' = '  +3.15%  '
    'D33' = '0.04670'
    'E33' = '  +3.04%  '
    'D34' = '2.662'
    'E34' = '  +1.85%  '
    'D35' = '0.9900'
    'E35' = '  +1.12%  '
    'D36' = '0.6199'
    'E36' = '  -0.32%  '
    'D37' = '2.688'
    'E37' = '  +0.15%  '
    'D38' = '0.01600'
    'E38' = '  +0.69%  '
    'D39' = '1.921'
    'E39' = '  -1.00%  '
    'D40' = '0.9996'
    'E40' = '  +0.05%  '
    'D41' = '99.81'
    'E41' = '  -0.44%  '
    'D42' = '0.3859'
    'E42' = '  -0.45%  '
    'D43' = '0.7469'
    'E43' = '  +2.02%  '
    'D44' = '5.003'
    'E44' = '  -0.13%  '
    'D45' = '0.1127'
    'E45' = '  +0.02%  '
    'D46' = '6.255'
    'E46' = '  -0.32%  '
    'D47' = '55.06'
    'E47' = '  +2.71%  '
    'D48' = '0.05213'
    'E48' = '  -2.51%  '
    'D49' = '30.61'
    'E49' = '  +1.48%  '
    'E50' = '  -2.12%  '
    'E51' = '  -0.50%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text format so values like "1.000" or "25.862.84" are not
    # reinterpreted as numbers/dates by Excel, matching the original inlineStr cells.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
